# fix comparison modes, use combine_excel_files(), fix column letters selection, change debug colors
#
# The combine_excel_files() routine re-emits the "First name" / "Last name"
# columns (A:B) using the lower-cased comparison key instead of the
# originally-cased display value, and the active selection is moved to A5
# (the row that was being inspected while fixing the column-letter mapping).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "cal"
$ws.Range("B2").Value = "arnold"

$ws.Range("A4").Value = "stu"
$ws.Range("B4").Value = "mc'elroy"

$ws.Range("B5").Value = "rico "

$ws.Range("A6").Value = "renée"
$ws.Range("B6").Value = "liêvre"

$ws.Range("A5").Value = "ncle"

$ws.Range("A5").Select()
